# The underlying survey records for rows 14-16 were re-sorted / re-numbered
# upstream (new observation IDs assigned), which rotates the data that was
# sitting in rows 14, 15 and 16 by one position:
#   new row 14 <- old row 15 (Kolflarnlav / Carbonicola anthracophila)
#   new row 15 <- old row 16 (Plattlummer, Antal=10)
#   new row 16 <- old row 14 (Plattlummer, Antal=5)
# Apply that as direct cell writes so every column lands on its new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 14 (becomes the old row 15 / "Kolflarnlav" record) ----
$ws.Range("A14").Value = 111380345
$ws.Range("B14").Value = 77267
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 6446
$ws.Range("F14").Value = "Kolflarnlav"
$ws.Range("G14").Value = "Carbonicola anthracophila"
$ws.Range("H14").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("I14").ClearContents()
$ws.Range("J14").ClearContents()
$ws.Range("Q14").Value = 364908.1256513004
$ws.Range("R14").Value = 6872135.474104149
$ws.Range("Z14").Value = "08:00"
$ws.Range("AB14").Value = "08:00"
$ws.Range("AC14").Value = "Växer på gammal kolad tallstubbe i kontinuitetsskog"
$ws.Range("AI14").Value = "Tallskog. Kontinuitetsskog"
$ws.Range("AJ14").Value = "tall"
$ws.Range("AK14").Value = "Pinus sylvestris"
$ws.Range("AO14").Value = "Pinus sylvestris"

# ---- Row 15 (becomes the old row 16 / "Plattlummer, Antal=10" record) ----
$ws.Range("A15").Value = 111379229
$ws.Range("B15").Value = 95538
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 221941
$ws.Range("F15").Value = "Plattlummer"
$ws.Range("G15").Value = "Lycopodium complanatum"
$ws.Range("H15").Value = "L."
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = "10"
$ws.Range("J15").NumberFormat = "@"
$ws.Range("J15").Value = "m²"
$ws.Range("Q15").Value = 364945.755472637
$ws.Range("R15").Value = 6872251.713583581
$ws.Range("AC15").Value = "Plattlummer växer i k-skog"
$ws.Range("AI15").Value = "Barrblandskog med gamla tallar och senvuxna granar. Kontinuitetsskog"
$ws.Range("AJ15").ClearContents()
$ws.Range("AK15").ClearContents()
$ws.Range("AO15").ClearContents()

# ---- Row 16 (becomes the old row 14 / "Plattlummer, Antal=5" record) ----
$ws.Range("A16").Value = 111379142
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "5"
$ws.Range("Q16").Value = 364964.1744805645
$ws.Range("R16").Value = 6872204.831332479
$ws.Range("Z16").Value = "00:00"
$ws.Range("AB16").Value = "00:00"
$ws.Range("AC16").Value = "Plattlummer växer runt gammal tall i k-skog"
